$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-09 Monday", 2)

$d.Content.Find.Execute("163÷7=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "674÷3=224, 2", 2)
$d.Content.Find.Execute("910÷3=303, 1", $true, $false, $false, $false, $false, $true, 1, $false, "218÷9=24, 2", 2)
$d.Content.Find.Execute("915÷3=305, 0", $true, $false, $false, $false, $false, $true, 1, $false, "188÷7=26, 6", 2)
$d.Content.Find.Execute("434÷2=217, 0", $true, $false, $false, $false, $false, $true, 1, $false, "426÷4=106, 2", 2)
$d.Content.Find.Execute("755÷5=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "624÷7=89, 1", 2)

$d.Content.Find.Execute("185÷9=20, 5", $true, $false, $false, $false, $false, $true, 1, $false, "312÷9=34, 6", 2)
$d.Content.Find.Execute("341÷9=37, 8", $true, $false, $false, $false, $false, $true, 1, $false, "116÷4=29, 0", 2)
$d.Content.Find.Execute("563÷6=93, 5", $true, $false, $false, $false, $false, $true, 1, $false, "438÷3=146, 0", 2)
$d.Content.Find.Execute("158÷4=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "595÷3=198, 1", 2)
$d.Content.Find.Execute("392÷3=130, 2", $true, $false, $false, $false, $false, $true, 1, $false, "644÷5=128, 4", 2)

$d.Content.Find.Execute("552÷2=276, 0", $true, $false, $false, $false, $false, $true, 1, $false, "720÷4=180, 0", 2)
$d.Content.Find.Execute("705÷7=100, 5", $true, $false, $false, $false, $false, $true, 1, $false, "569÷6=94, 5", 2)
$d.Content.Find.Execute("431÷4=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "313÷5=62, 3", 2)
$d.Content.Find.Execute("523÷6=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "986÷5=197, 1", 2)
$d.Content.Find.Execute("963÷3=321, 0", $true, $false, $false, $false, $false, $true, 1, $false, "714÷9=79, 3", 2)

$d.Content.Find.Execute("174÷5=34, 4", $true, $false, $false, $false, $false, $true, 1, $false, "874÷2=437, 0", 2)
$d.Content.Find.Execute("482÷9=53, 5", $true, $false, $false, $false, $false, $true, 1, $false, "945÷5=189, 0", 2)
$d.Content.Find.Execute("753÷4=188, 1", $true, $false, $false, $false, $false, $true, 1, $false, "637÷9=70, 7", 2)
$d.Content.Find.Execute("763÷4=190, 3", $true, $false, $false, $false, $false, $true, 1, $false, "793÷3=264, 1", 2)
$d.Content.Find.Execute("690÷3=230, 0", $true, $false, $false, $false, $false, $true, 1, $false, "156÷6=26, 0", 2)

$d.Content.Find.Execute("297÷2=148, 1", $true, $false, $false, $false, $false, $true, 1, $false, "373÷8=46, 5", 2)
$d.Content.Find.Execute("952÷7=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "110÷2=55, 0", 2)
$d.Content.Find.Execute("774÷6=129, 0", $true, $false, $false, $false, $false, $true, 1, $false, "159÷4=39, 3", 2)
$d.Content.Find.Execute("320÷7=45, 5", $true, $false, $false, $false, $false, $true, 1, $false, "764÷7=109, 1", 2)
$d.Content.Find.Execute("710÷9=78, 8", $true, $false, $false, $false, $false, $true, 1, $false, "881÷5=176, 1", 2)
